$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objetivosPt = 'A disciplina tem o objetivo de apresentar ao estudante informações a respeito das características dos solos, particularmente os existentes na região tropical, e dos fenômenos físicos que nele ocorrem, a fim de capacitá-lo a compreender a importância dos fatores pedológicos, físicos e hídricos na preservação do ambiente.'
$docentesValue = '5840942 - Marco Aurélio Kondracki de Alcântara'
$programaResumidoPt = 'Introdução. Formação do solo. Atributos físicos do solo. Classificação do solo. Água do solo. Aula Prática: Descrição de perfil no campo. Aula Prática: Caracterização e métodos de determinação de atributos físicos e hídricos do solo.'
$programaPt = 'INTRODUÇÃO. Conceitos Básicos. O perfil de solo. Definição e notação de horizontes e camadas. FORMAÇÃO DO SOLO. Fatores e processos de formação. Intemperismo. ATRIBUTOS FÍSICOS DO SOLO. Composição volumétrica, granulometria e textura, estrutura e agregação, cor, porosidade, densidade e compactação, consistência. CLASSIFICAÇÃO DO SOLO. Sistema brasileiro de classificação de solos. Principais atributos morfológicos. Principais Classes de Solos. ÁGUA DO SOLO. Conceito e importância. Constantes de umidade. Potencial total da água do solo e seus componentes. Curva característica da água do solo. Movimento da Água e de solutos no Solo. Aula prática de campo: Descrição de perfil no campo. Aula prática de laboratório: Caracterização e métodos de determinação de atributos físicos e hídricos do solo.'
$metodoText = 'A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas.'
$criterioText = 'O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas.'
$normaText = 'Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.'
$biblioText = 'Bibliografia básica:1. MEURER, E.J. Fundamentos de Química do Solo, 3.ed. Porto Alegre: Editora Evangraf, 2010. 264p.2. ERNANI, P.R. Química do Solo e Disponibilidade de Nutrientes. Lages: Udesc, 1. ed. , 2008. v.1. 230 p.Bibliografia complementar:1. CAMARGO, O.A. de; MONIZ, A.C.; JORGE, J.A.; VALADARES, J.M.A.S. Métodos de analise química, mineralógica e física de solos do Instituto Agronômico de Campinas. Campinas, Instituto Agronômico, 2009. 77 p. (Boletim técnico, 106, Edição revista e atualizada).2. EMBRAPA. Centro Nacional de Pesquisa de Solos (Rio de Janeiro, RJ). Sistema Brasileiro de Classificação de Solos. Brasília: Embrapa Produção da Informação; Rio de Janeiro: Embrapa Solos, 2006. 306p.3. MELO, V.F.; ALLEONI, L.R.F. (Eds.). Química e mineralogia do solo. v.1: Conceitos básicos. Viçosa: SBCS, 2009. 595p. 5. MELO, V.F.; ALLEONI, L.R.F. (Eds.). Química e mineralogia do solo. v.2: Aplicações. Viçosa: SBCS, 2009. 685p.4. NOVAIS, R.F.; ALVAREZ V., V.H.; BARROS, N.F.; FONTES, R.L.F.; CANTARUTTI,R.B.; NEVES, J.C.L. Fertilidade do Solo. Visconde do Rio Branco: Gráfica Suprema, 2007. 1017p.5. QUAGGIO, J. A. Acidez e calagem em solos tropicais. Instituto Agronômico. 111p. (2000).6. RAIJ, B. van; ANDRADE, J.C. de; CANTARELLA, H.; QUAGGIO, J.A. Análise química para avaliação da fertilidade de solos tropicais. Raij, B. van, Andrade, J.C. de, Cantarella, H. e Quaggio, J.A. (ed.). Campinas, Instituto Agronômico, 2001. 285p.7. SANTOS, G.A; SILVA, L.S.; CANELLAS, L.P.; CAMARGO, F.A.O. (Eds). Fundamentos da matéria orgânica do solo: ecossistemas tropicais e subtropicais. Porto Alegre: Genesis. 2a Edição. 2008. 636p.'


# Insert a new row at position 13. This shifts the existing rows 13-21 down to 14-22,
# which both fixes the pre-existing label/value misalignment and matches the row-height
# pattern of the target layout (ht=60/120 flags move down together with their rows).
$ws.Rows("13:13").Insert()

# The newly inserted row 13 picked up a stray A13 cell (copied format from the old A13).
# Target layout has no A13 - row 13 only carries the "Docentes responsaveis" value in B/C.
$ws.Range("A13").Clear()

# Row 10: Objetivos: - fill in the correct Portuguese objectives paragraph (was blank/wrong).
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# Row 13 (new): value-only row holding the "Docentes responsaveis:" (row 12 label) answer.
# B13/C13 are brand-new cells (the inserted row only seeded A13), so pull in the normal
# column formatting (styles 2/3) from a template row before writing the value.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = $docentesValue
$ws.Range("C13").Value = $docentesValue

# Row 14: Programa resumido: - replace placeholder "Semestral" with the real summary text.
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt

# Row 16: Programa: - replace the wrong date value with the full Portuguese syllabus text.
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# Row 19: Metodo: - correct value (was wrongly carrying the "Docentes responsaveis" text).
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Row 20: Criterio: - correct value (previously held the Metodo text).
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Row 21: Norma de recuperacao: - correct value (previously held the Criterio text).
$ws.Range("B21").Value = $normaText
$ws.Range("C21").Value = $normaText

# Row 22: Bibliografia: - correct value (previously held the Norma de recuperacao text).
$ws.Range("B22").Value = $biblioText
$ws.Range("C22").Value = $biblioText
